$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so values like "19.40" or
# "36.654.19" are stored as literal text (matching the source data) instead
# of being auto-coerced to numbers/dates by Excel's smart input parsing.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '36.654.19'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '1.957.96'
$ws.Range("E3").Value = '  +0.90%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '244.47'
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("D7").Value = '61.61'
$ws.Range("E7").Value = '  +8.13%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +4.84%  '
$ws.Range("D10").Value = '0.0793'
$ws.Range("E10").Value = '  -6.39%  '
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").Value = '14.24'
$ws.Range("E12").Value = '  +6.20%  '
$ws.Range("D13").Value = '21.94'
$ws.Range("E13").Value = '  +3.22%  '
$ws.Range("E14").Value = '  +2.92%  '
$ws.Range("D15").Value = '2.229.98'
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("E16").Value = '  +3.01%  '
$ws.Range("D17").Value = '1.952.56'
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("D18").Value = '36.537.14'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").Value = '69.71'
$ws.Range("E19").Value = '  +0.89%  '
$ws.Range("D20").Value = '0.0₃0852'
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").Value = '230.01'
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("E22").Value = '  +1.93%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").Value = '  +5.07%  '
$ws.Range("E25").Value = '  +2.95%  '
$ws.Range("E26").Value = '  +7.72%  '
$ws.Range("D27").Value = '9.17'
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").Value = '160.48'
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").Value = '19.40'
$ws.Range("E29").Value = '  +1.03%  '
$ws.Range("D30").Value = '1.29'
$ws.Range("E30").Value = '  +17.64%  '
$ws.Range("E31").Value = '  +1.22%  '
$ws.Range("D32").Value = '4.76'
$ws.Range("E32").Value = '  +4.61%  '
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("E34").Value = '  +7.05%  '
$ws.Range("D35").Value = '3.55'
$ws.Range("E35").Value = '  +13.56%  '
$ws.Range("E36").Value = '  +4.76%  '
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").Value = '1.78'
$ws.Range("E38").Value = '  -0.91%  '
$ws.Range("D39").Value = '5.51'
$ws.Range("E39").Value = '  -9.34%  '
$ws.Range("D40").Value = '0.0979'
$ws.Range("E40").Value = '  -0.93%  '
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("E42").Value = '  +2.09%  '
$ws.Range("E43").Value = '  +0.95%  '
$ws.Range("D44").Value = '15.99'
$ws.Range("E44").Value = '  +2.64%  '
$ws.Range("D45").Value = '1.366.53'
$ws.Range("E45").Value = '  +2.02%  '
$ws.Range("D46").Value = '88.52'
$ws.Range("E46").Value = '  +2.58%  '
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").Value = '45.55'
$ws.Range("E50").Value = '  +5.87%  '
$ws.Range("D51").Value = '2.127.91'
$ws.Range("E51").Value = '  +0.57%  '

# Restore the original (unstyled) formatting now that the text values are set.
$ws.Range("D2:E51").ClearFormats()
